# Apply "Added Code For inventory Transactions" edit:
#  - remove the value from I2 (was "Pro-1")
#  - duplicate row 12 (as originally authored) into a new row 14
#  - change I12's value to the new string "MNK60"
#  - move the active selection to row 12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot row 12's original contents via .Value2 (the only reliably
# type-preserving getter here) before row 12 itself is modified.
$r12A = $ws.Range("A12").Value2
$r12B = $ws.Range("B12").Value2
$r12C = $ws.Range("C12").Value2
$r12D = $ws.Range("D12").Value2
$r12E = $ws.Range("E12").Value2
$r12F = $ws.Range("F12").Value2
$r12G = $ws.Range("G12").Value2
$r12H = $ws.Range("H12").Value2
$r12I = $ws.Range("I12").Value2
$r12J = $ws.Range("J12").Value2
$r12K = $ws.Range("K12").Value2
$r12L = $ws.Range("L12").Value2
$r12M = $ws.Range("M12").Value2

# New row 14 = a copy of row 12's original content
$ws.Range("A14").Value2 = $r12A
$ws.Range("B14").Value2 = $r12B
$ws.Range("C14").Value2 = $r12C
$ws.Range("D14").Value2 = $r12D
$ws.Range("E14").Value2 = $r12E
$ws.Range("F14").Value2 = $r12F
$ws.Range("G14").Value2 = $r12G
$ws.Range("H14").Value2 = $r12H
$ws.Range("I14").Value2 = $r12I
$ws.Range("J14").Value2 = $r12J
$ws.Range("K14").Value2 = $r12K
$ws.Range("L14").Value2 = $r12L
$ws.Range("M14").Value2 = $r12M

# K12's cell style (the date number format) needs to carry over to K14 too;
# copy formats only so no duplicate style entries get created and K14's
# value (set above) is left untouched.
$ws.Range("K12").Copy()
$ws.Range("K14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Remove the value that used to live in I2 ("Pro-1")
$ws.Range("I2").ClearContents()

# Row 12's I cell now references the newly introduced shared string "MNK60"
$ws.Range("I12").Value2 = "MNK60"

# Reflect the new selection state stored in the workbook (active cell A12,
# with the whole row selected)
$ws.Range("A12:XFD12").Select()
